$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13542.857
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 13542.857
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 13542.857
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -14478.857
$ws.Range("H23").Value = 13542.857
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 13542.857
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 13542.857
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -14010.857
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 3000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -3562
$ws.Range("H86").Value = 11348.75
$ws.Range("I86").Value = 2700
$ws.Range("J86").Value = 19997.5
$ws.Range("K86").Value = 2700
$ws.Range("L86").Value = 19997.5
$ws.Range("M86").Value = -1577
$ws.Range("N86").Value = -22243.5
$ws.Range("H89").Value = 11348.75
$ws.Range("I89").Value = 2700
$ws.Range("J89").Value = 19997.5
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 99987.5
$ws.Range("M89").Value = -7884
$ws.Range("N89").Value = -111219.5
$ws.Range("H121").Value = 1912
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 1890
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 5670
$ws.Range("M121").Value = -4253
$ws.Range("N121").Value = -9164
$ws.Range("H135").Value = 1234.1515
$ws.Range("I135").Value = 636.4375
$ws.Range("J135").Value = 1796.7059
$ws.Range("K135").Value = 5727.9375
$ws.Range("L135").Value = 16170.3531
$ws.Range("M135").Value = -3192.9375
$ws.Range("N135").Value = -21240.3531
$ws.Range("H138").Value = 1704.4
$ws.Range("J138").Value = 1938.2174
$ws.Range("L138").Value = 5814.6522
$ws.Range("N138").Value = -16094.6522
$ws.Range("H140").Value = 60996.668
$ws.Range("J140").Value = 60996.668
$ws.Range("L140").Value = 60996.668
$ws.Range("N140").Value = -71356.66800000001
$ws.Range("H141").Value = 4034.8125
$ws.Range("I141").Value = 4484.625
$ws.Range("J141").Value = 3585
$ws.Range("K141").Value = 13453.875
$ws.Range("L141").Value = 10755
$ws.Range("M141").Value = -8273.875
$ws.Range("N141").Value = -21115

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 165
$ws.Range("I4").Value = 90.5
$ws.Range("K4").Value = 90.5
$ws.Range("M4").Value = 25.5
$ws.Range("H32").Value = 23920.87
$ws.Range("I32").Value = 5277.84
$ws.Range("J32").Value = 179279.44
$ws.Range("K32").Value = 5277.84
$ws.Range("L32").Value = 179279.44
$ws.Range("M32").Value = -4990.84
$ws.Range("N32").Value = -179853.44
$ws.Range("H133").Value = 37000
$ws.Range("J133").Value = 37000
$ws.Range("L133").Value = 37000
$ws.Range("N133").Value = -42060
$ws.Range("H135").Value = 35054.9
$ws.Range("J135").Value = 35054.9
$ws.Range("L135").Value = 35054.9
$ws.Range("N135").Value = -45194.9

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2769.375
$ws.Range("I54").Value = 2352.2856
$ws.Range("K54").Value = 2352.2856
$ws.Range("M54").Value = -1868.2856
$ws.Range("H105").Value = 119620.06
$ws.Range("I105").Value = 79224.53999999999
$ws.Range("K105").Value = 79224.53999999999
$ws.Range("M105").Value = -77477.53999999999
$ws.Range("H106").Value = 31500
$ws.Range("J106").Value = 31500
$ws.Range("L106").Value = 31500
$ws.Range("N106").Value = -34024
$ws.Range("H134").Value = 2999.54
$ws.Range("I134").Value = 2492.682
$ws.Range("J134").Value = 6716.5
$ws.Range("K134").Value = 7478.045999999999
$ws.Range("L134").Value = 20149.5
$ws.Range("M134").Value = -4943.045999999999
$ws.Range("N134").Value = -25219.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1931.1428
$ws.Range("I23").Value = 1379.5
$ws.Range("J23").Value = 2666.6667
$ws.Range("K23").Value = 1379.5
$ws.Range("L23").Value = 2666.6667
$ws.Range("M23").Value = -1139.5
$ws.Range("N23").Value = -3146.6667
$ws.Range("H27").Value = 1931.1428
$ws.Range("I27").Value = 1379.5
$ws.Range("J27").Value = 2666.6667
$ws.Range("K27").Value = 1379.5
$ws.Range("L27").Value = 2666.6667
$ws.Range("M27").Value = -1187.5
$ws.Range("N27").Value = -3050.6667
$ws.Range("H105").Value = 1514.1
$ws.Range("I105").Value = 1696
$ws.Range("J105").Value = 1332.2
$ws.Range("K105").Value = 1696
$ws.Range("L105").Value = 1332.2
$ws.Range("M105").Value = 51
$ws.Range("N105").Value = -4826.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2217.2144
$ws.Range("I33").Value = 1699.8334
$ws.Range("J33").Value = 2605.25
$ws.Range("K33").Value = 10199.0004
$ws.Range("L33").Value = 15631.5
$ws.Range("M33").Value = -9916.000400000001
$ws.Range("N33").Value = -16197.5
$ws.Range("H132").Value = 1256.3572
$ws.Range("I132").Value = 580.36365
$ws.Range("J132").Value = 3735
$ws.Range("K132").Value = 5223.27285
$ws.Range("L132").Value = 33615
$ws.Range("M132").Value = -2693.27285
$ws.Range("N132").Value = -38675

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6104.8
$ws.Range("I99").Value = 2836.5557
$ws.Range("J99").Value = 35519
$ws.Range("K99").Value = 2836.5557
$ws.Range("L99").Value = 35519
$ws.Range("M99").Value = -590.5556999999999
$ws.Range("N99").Value = -40011
$ws.Range("H107").Value = 631514.9
$ws.Range("J107").Value = 1262737
$ws.Range("L107").Value = 1262737
$ws.Range("N107").Value = -1266577
$ws.Range("H135").Value = 32122
$ws.Range("J135").Value = 32122
$ws.Range("L135").Value = 32122
$ws.Range("N135").Value = -42262

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3633.3333
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 4950
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 4950
$ws.Range("M9").Value = -776
$ws.Range("N9").Value = -5398
$ws.Range("H22").Value = 1074.5128
$ws.Range("I22").Value = 1821.8
$ws.Range("J22").Value = 816.8276
$ws.Range("K22").Value = 1821.8
$ws.Range("L22").Value = 816.8276
$ws.Range("M22").Value = -1526.8
$ws.Range("N22").Value = -1406.8276
$ws.Range("H27").Value = 1074.5128
$ws.Range("I27").Value = 1821.8
$ws.Range("J27").Value = 816.8276
$ws.Range("K27").Value = 1821.8
$ws.Range("L27").Value = 816.8276
$ws.Range("M27").Value = -1714.8
$ws.Range("N27").Value = -1030.8276
$ws.Range("H38").Value = 14133.333
$ws.Range("I38").Value = 8000
$ws.Range("K38").Value = 8000
$ws.Range("M38").Value = -7590
$ws.Range("H55").Value = 1000.03705
$ws.Range("I55").Value = 1673.1
$ws.Range("J55").Value = 604.1177
$ws.Range("K55").Value = 1673.1
$ws.Range("L55").Value = 604.1177
$ws.Range("M55").Value = -1500.1
$ws.Range("N55").Value = -950.1177
$ws.Range("H74").Value = 19299.2
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 19299.2
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 19299.2
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -21295.2
$ws.Range("H77").Value = 19299.2
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 19299.2
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 57897.60000000001
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -67881.60000000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1668.4231
$ws.Range("I122").Value = 1195.25
$ws.Range("J122").Value = 2425.5
$ws.Range("K122").Value = 3585.75
$ws.Range("L122").Value = 7276.5
$ws.Range("M122").Value = -1135.75
$ws.Range("N122").Value = -12176.5
